$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the authoritative XML diff.
# Values are written as TEXT (matching the source inlineStr cells) by forcing a
# "@" (text) number format before the assignment, then clearing the format again
# so the cell keeps using the sheet default style (no explicit style index).
$changes = [ordered]@{
    'D2' = '230.96'
    'D3' = '22.43'
    'D4' = '5.515'
    'D5' = '0.05556'
    'D6' = '3.396'
    'D7' = '6.491'
    'D9' = '0.7899'
    'B10' = 'One'
    'C10' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D10' = '0.01149'
    'E10' = '9OneONEBestin24h'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1395'
    'E11' = '10WazirXWRX'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.07454'
    'E12' = '11MandalaExchangeTokenMDX'
    'B13' = 'LiechtensteinCryptoassetsExchange'
    'C13' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D13' = '0.03135'
    'E13' = '12LiechtensteinCryptoassetsExchangeLCX'
    'B14' = 'BitrueCoin'
    'C14' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D14' = '0.02912'
    'E14' = '13BitrueCoinBTR'
    'B15' = 'BitMartToken'
    'C15' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D15' = '0.09276'
    'E15' = '14BitMartTokenBMX'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '0.001666'
    'E16' = '15BitForexTokenBF'
    'B17' = 'MCDex'
    'C17' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D17' = '3.277'
    'E17' = '16MCDexMCB'
    'B18' = 'CoinExToken'
    'C18' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D18' = '0.04743'
    'E18' = '17CoinExTokenCET'
    'D19' = '0.006249'
    'D20' = '0.005249'
    'D22' = '0.0001507'
    'D23' = '3.695'
    'D26' = '0.1290'
    'D27' = '0.0006161'
    'E27' = '26UpBotsUBXT'
    'D41' = '0.007152'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D42' = '0.003349'
    'E42' = '41CEJICEJI'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D43' = '0.1026'
    'E43' = '42BKEXTokenBKK'
    'D44' = '0.009131'
    'E44' = '43LocalTradersLCT'
    'E45' = '44ACDXExchangeACXTWorstin24h'
    'D46' = '0.00005519'
    'D47' = '0.00000000753'
    'D48' = '0.6782'
    'D49' = '0.09415'
    'D50' = '0.00002109'
    'D51' = '0.01014'
}

foreach ($addr in $changes.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $changes[$addr]
    $rng.ClearFormats()
}
